$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Locate the existing hyperlink that sits on M2 (keeps the same r:id / target
# address - only the cell text + displayed caption are changing).
$m2Link = $null
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 13) {
        $m2Link = $hl
    }
}

# Update the cell text shown in M2 (was "https://mirandakate.cabitest21.com",
# now "https://mirandakate.cabitest5.com"). The hyperlink target itself is
# left untouched.
$ws.Range("M2").Value = "https://mirandakate.cabitest5.com"

# Clear the hyperlink's explicit display caption so it no longer diverges
# from the cell text (drops the now-stale display="..." attribute).
if ($m2Link -ne $null) {
    $m2Link.TextToDisplay = ""
}

# Move the active selection to M2, matching where the edit was made.
$ws.Range("M2").Select()
